# Adds a new "configMode" entry to the functions list (sheet "functions"),
# and wires it up as the selected function for row 31 on the "midimap"
# sheet (J31 = configMode, K31 = 2), matching commit "added config mode selector".

$wb = $excel.ActiveWorkbook

$wsMidimap   = $wb.Worksheets.Item("midimap")
$wsFunctions = $wb.Worksheets.Item("functions")

# 1) Insert the new function name alphabetically into the "functions" sheet.
#    Existing sorted list: agc, cwxLoopBtnL, cwxLoopBtnR, fadePan, filters, ...
#    "configMode" sorts right before "cwxLoopBtnL", i.e. at row 3.
$wsFunctions.Rows.Item(3).Insert()
$wsFunctions.Range("A3").Value = "configMode"

# 2) On the "midimap" sheet, assign the new function to row 31: set J31 to
#    the new function name and bump K31 (its mode/value column) from 0 to 2.
$wsMidimap.Range("J31").Value = "configMode"
$wsMidimap.Range("K31").Value = 2

# 3) Restore view state: active sheet/selection on "midimap" is K31, and the
#    "functions" sheet's own selection moves to B46.
$wsFunctions.Range("B46").Select() | Out-Null
$wsMidimap.Activate() | Out-Null
$wsMidimap.Range("K31").Select() | Out-Null
